# Updated cryptos list on Mon Aug 21 10:38:23 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of the
# crypto table on Sheet1 with the latest scrape. Two coins (Cronos /
# EnergySwap, rows 48-49) also swapped rank order, so their Coin name, Link,
# Price and Volume cells are all rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few Price values need to stay exactly as scraped (trailing zeros such as
# "99.90" / "3.530", or long decimals like "0.00000000109") instead of being
# silently normalised into numbers by Excel's type inference. Pre-formatting
# just those cells as Text keeps the literal string intact.
$textCells = @("D11", "D32", "D43", "D45", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.113.16"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.675.38"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
$ws.Range("D5").Value = "212.43"
$ws.Range("E5").Value = "  -2.68%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.5266"
$ws.Range("E6").Value = "  -4.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.14%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").Value = "  -1.15%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06311"
$ws.Range("E9").Value = "  -2.85%  "

# Row 10 - Solana
$ws.Range("D10").Value = "21.25"
$ws.Range("E10").Value = "  -3.85%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07610"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.694.82"
$ws.Range("E12").Value = "  +0.73%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "4.501"
$ws.Range("E13").Value = "  -1.02%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.5679"
$ws.Range("E14").Value = "  -2.30%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.000008163"
$ws.Range("E15").Value = "  -3.50%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.71"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.133.34"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  -0.14%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "4.832"
$ws.Range("E19").Value = "  -2.12%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "10.61"
$ws.Range("E20").Value = "  -2.85%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "188.92"
$ws.Range("E21").Value = "  -1.20%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.191"
$ws.Range("E22").Value = "  -0.69%  "

# Row 23 - BinanceUSD
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - Monero
$ws.Range("D24").Value = "148.77"
$ws.Range("E24").Value = "  +0.90%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "0.1251"
$ws.Range("E25").Value = "  -5.20%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "7.628"
$ws.Range("E26").Value = "  -3.40%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - Hedera
$ws.Range("D28").Value = "0.06378"
$ws.Range("E28").Value = "  +0.62%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "1.358"
$ws.Range("E29").Value = "  -2.43%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "1.296"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "3.538"
$ws.Range("E31").Value = "  -1.55%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "1.659"
$ws.Range("E33").Value = "  -0.46%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.011"
$ws.Range("E34").Value = "  -2.63%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.419"
$ws.Range("E35").Value = "  +0.73%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.6018"
$ws.Range("E36").Value = "  -3.08%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38 - FraxShare
$ws.Range("D38").Value = "6.134"
$ws.Range("E38").Value = "  -1.64%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.01624"
$ws.Range("E39").Value = "  -0.03%  "

# Row 40 - Maker
$ws.Range("D40").Value = "1.089.01"
$ws.Range("E40").Value = "  -2.08%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.8714"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - PaxDollar
$ws.Range("D42").Value = "1.006"
$ws.Range("E42").Value = "  -0.99%  "

# Row 43 - Quant
$ws.Range("D43").Value = "99.90"
$ws.Range("E43").Value = "  -0.83%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.829.10"
$ws.Range("E44").Value = "  -0.14%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  +0.71%  "

# Row 46 - Aave
$ws.Range("D46").Value = "56.97"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47 - Frax
$ws.Range("E47").Value = "  -0.05%  "

# Row 48 - was Cronos, now EnergySwap (rank swap with row 49)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.007"
$ws.Range("E48").Value = "  -2.13%  "

# Row 49 - was EnergySwap, now Cronos (rank swap with row 48)
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.05253"
$ws.Range("E49").Value = "  -0.38%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "0.4275"
$ws.Range("E50").Value = "  -0.42%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "5.930"
$ws.Range("E51").Value = "  -2.36%  "
